# Refactor & code for card view
# Add a new "maxGroupNum" (最大卡组数量) config row to the "main" sheet
# and update the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

# Populate row 11 with the new config entry: name / key / value
$ws.Cells.Item(11, 1).Value = "最大卡组数量"
$ws.Cells.Item(11, 2).Value = "maxGroupNum"
$ws.Cells.Item(11, 3).Value = 18

# Update the selected/active cell on the sheet to D12
$ws.Range("D12").Select()
